# Auto-generated Excel COM-interop script to apply the "Updated symbol list" diff
# (cryptos.xlsx) — updates Coin/Link (B/C) and Price/Volume(1h) (D/E) cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column -> new text value, taken from the target diff.
$changes = @{
    2 = @{ "D" = "294.59"; "E" = "1.03%" }
    3 = @{ "D" = "40.06"; "E" = "0.15%" }
    4 = @{ "D" = "5.018"; "E" = "-0.04%" }
    5 = @{ "D" = "0.07355"; "E" = "-0.01%" }
    6 = @{ "B" = "GateToken"; "C" = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"; "D" = "4.290"; "E" = "0.05%" }
    7 = @{ "B" = "FTXToken"; "C" = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"; "D" = "1.546"; "E" = "-0.69%" }
    8 = @{ "B" = "MXToken"; "C" = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"; "D" = "0.9237"; "E" = "0.56%" }
    9 = @{ "B" = "BTSEToken"; "C" = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"; "D" = "2.382"; "E" = "-0.63%" }
    10 = @{ "B" = "LiechtensteinCryptoassetsExchange"; "C" = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"; "D" = "0.1180"; "E" = "-0.76%" }
    11 = @{ "B" = "WazirX"; "C" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"; "D" = "0.1784"; "E" = "2.57%" }
    12 = @{ "B" = "MandalaExchangeToken"; "C" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"; "D" = "0.08680"; "E" = "-0.60%" }
    13 = @{ "B" = "BitrueCoin"; "C" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"; "D" = "0.04237"; "E" = "1.48%" }
    14 = @{ "B" = "BitMartToken"; "C" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"; "D" = "0.1054"; "E" = "0.13%" }
    15 = @{ "B" = "BitForexToken"; "C" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"; "D" = "0.001254"; "E" = "-1.79%" }
    16 = @{ "B" = "TigerCash"; "C" = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"; "D" = "0.005980"; "E" = "3.52%" }
    17 = @{ "B" = "LEO"; "C" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"; "D" = "3.367"; "E" = "-1.15%" }
    18 = @{ "D" = "0.3268"; "E" = "-0.86%" }
    19 = @{ "D" = "7.867"; "E" = "3.64%" }
    20 = @{ "D" = "0.1383"; "E" = "2.23%" }
    22 = @{ "D" = "0.03938"; "E" = "2.31%" }
    23 = @{ "D" = "0.001268"; "E" = "-0.66%" }
    24 = @{ "D" = "0.003816"; "E" = "-1.91%" }
    25 = @{ "D" = "0.0001232"; "E" = "-4.00%" }
    26 = @{ "E" = "0.01%" }
    38 = @{ "D" = "0.02339"; "E" = "0.18%" }
    39 = @{ "D" = "0.05078"; "E" = "1.09%" }
    40 = @{ "D" = "0.006139"; "E" = "20.15%" }
    41 = @{ "D" = "0.007783"; "E" = "1.29%" }
    42 = @{ "D" = "0.1290"; "E" = "1.67%" }
    43 = @{ "D" = "0.007399"; "E" = "0.08%" }
    44 = @{ "D" = "0.007185"; "E" = "-6.74%" }
    46 = @{ "D" = "0.00006128"; "E" = "-6.07%" }
    47 = @{ "E" = "0.02%" }
    48 = @{ "D" = "0.04659"; "E" = "-81.50%" }
    49 = @{ "E" = "0.02%" }
    50 = @{ "D" = "0.00002104"; "E" = "0.02%" }
    51 = @{ "E" = "0.02%" }
}

# Columns B/C hold plain text (names/links); D/E hold numeric-looking text
# (price/volume) that must stay stored as text, matching the workbook's
# original inline-string cells, so force a Text number format before writing.
$textNumberCols = @("D", "E")

foreach ($rowNum in ($changes.Keys | Sort-Object {[int]$_})) {
    $rowChanges = $changes[$rowNum]
    foreach ($col in $rowChanges.Keys) {
        $cell = $ws.Range("$col$rowNum")
        if ($textNumberCols -contains $col) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $rowChanges[$col]
    }
}
